# Insert a new "total_concentrations" worksheet right after
# "input_concentrations" (and before "equilibrium_concentrations"),
# populate it with the total-concentration inputs, and restore the
# original active-sheet selection.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("input_concentrations")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "total_concentrations"

# Header row
$ws.Range("A1").Value = "H"
$ws.Range("B1").Value = "PO4"
$ws.Range("C1").Value = "Cu"

# Data rows
$ws.Range("A2").Value = 0.01
$ws.Range("B2").Value = 0.01
$ws.Range("C2").Value = 0.01

$ws.Range("A3").Value = 0.02
$ws.Range("B3").Value = 0.01
$ws.Range("C3").Value = 0.01

$ws.Range("A4").Value = 0.03
$ws.Range("B4").Value = 0.01
$ws.Range("C4").Value = 0.01

# Adding the sheet shifts the active tab onto the new sheet; restore the
# workbook's original selection (first sheet was active before the edit).
$wb.Worksheets.Item("input_stoich_coefficients").Activate()
